# Improve the email's format following Varrun Gupta's suggestion:
# append four more "Alex Chi / 100%" rows to the "21-09-28" dogfooding sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("21-09-28")

# Copy the existing "Alex Chi / 100%" row (row 3) down into the four new
# rows so the new cells inherit the same (default) formatting/type as the
# pre-existing rows, rather than being re-interpreted (e.g. "100%" turning
# into a numeric percentage).
$ws.Range("A3:B3").Copy($ws.Range("A4:B4"))
$ws.Range("A3:B3").Copy($ws.Range("A5:B5"))
$ws.Range("A3:B3").Copy($ws.Range("A6:B6"))
$ws.Range("A3:B3").Copy($ws.Range("A7:B7"))
